# "Fix indentation in paper template"
#
# 1) Move the (collapsed) "_GoBack" bookmark so it sits after the
#    "Hello world." run instead of before it.
# 2) BodyText style: change the paragraph indent from a *left* indent
#    of 284 twips (14.2pt) to a *first-line* indent of 284 twips.
# 3) FirstParagraph style (based on BodyText): explicitly zero out the
#    first-line indent it now inherits from BodyText.

$d = $word.ActiveDocument

# --- 1. Reposition the _GoBack bookmark -----------------------------
# A bookmark collapsed exactly at "end of paragraph text" gets pinned to
# the paragraph boundaries when written back out, so we can't just grab
# $d.Range(12,12) and re-add it there directly. Work around it by
# temporarily appending a sentinel character, anchoring the collapsed
# bookmark right before that sentinel (i.e. right after "Hello world."),
# and then removing the sentinel again.
$body = $d.Paragraphs(1).Range
$endOfText = $body.End - 1   # just before the paragraph mark

$d.Bookmarks("_GoBack").Delete()

$sentinelStart = $endOfText
$d.Range($sentinelStart, $sentinelStart).InsertAfter("@")

$bmRange = $d.Range($endOfText, $endOfText)
$d.Bookmarks.Add("_GoBack", $bmRange)

$d.Range($endOfText, $endOfText + 1).Delete()

# --- 2. BodyText: left indent -> first-line indent -------------------
$bodyText = $d.Styles("Body Text")
$bodyText.ParagraphFormat.LeftIndent = 0
$bodyText.ParagraphFormat.FirstLineIndent = 14.2

# --- 3. FirstParagraph: cancel the inherited first-line indent -------
$firstParagraph = $d.Styles("First Paragraph")
$firstParagraph.ParagraphFormat.FirstLineIndent = 0
